# Automatic update of files.
# - Column C ("Förändrad") advances by one day (46063 -> 46064) on every data row.
# - Rows 7-11 ("Beteckning"/A, "Datum"/B, "Area (ha)"/G) rotate up by one: the
#   record that was in row 7 moves down to row 11, and the records that were
#   in rows 8-11 each shift up into rows 7-10.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bump the "Förändrad" date (column C) for every data row (2-11) by one day.
for ($r = 2; $r -le 11; $r++) {
    $ws.Cells.Item($r, 3).Value = 46064
}

# Rotate the record data (columns A, B, G) for rows 7-11 up by one position,
# wrapping row 7's original record around to row 11.
$beteckning = @("A 33037-2025", "A 33033-2025", "A 6314-2022", "A 25610-2024", "A 57810-2022")
$datum      = @(45840.39623842593, 45840.39188657407, 44600, 45463, 44897)
$area       = @(0.8, 0.7, 3, 2.9, 3.3)

for ($i = 0; $i -lt 5; $i++) {
    $r = 7 + $i
    $ws.Cells.Item($r, 1).Value = $beteckning[$i]
    $ws.Cells.Item($r, 2).Value = $datum[$i]
    $ws.Cells.Item($r, 7).Value = $area[$i]
}
